$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (RF) updated values
$ws.Range("B3").Value = 0.318
$ws.Range("C3").Value = 0.016
$ws.Range("D3").Value = 0.433
$ws.Range("E3").Value = 0.658
$ws.Range("F3").Value = 0.766
$ws.Range("G3").Value = 0.591

# Row 4 (NN) updated values
$ws.Range("B4").Value = -0.04
$ws.Range("C4").Value = -0.5
$ws.Range("D4").Value = 0.66
$ws.Range("E4").Value = 0.812
$ws.Range("F4").Value = 0.848
$ws.Range("G4").Value = 0.455

# Row 5 (RNN) updated values
$ws.Range("B5").Value = -0.046
$ws.Range("C5").Value = -0.292
$ws.Range("D5").Value = 0.597
$ws.Range("E5").Value = 0.773
$ws.Range("F5").Value = 0.753
$ws.Range("G5").Value = 0.372
